# "added the ability to edit the left table on leftclick"
#
# Append 5 new transaction rows (11-15) to the left table on Sheet1.
# Category | Name | Date | Price | Account
#
# A few of the values look like numbers or dates ("1202.0", "2023-03-14",
# etc.) but in the source workbook they are stored as plain text (shared
# strings), matching the style of every other row already on the sheet.
# Excel's COM layer auto-detects such literals as numbers/dates, so for
# those cells we briefly force a text format before assigning the value,
# then clear the cell formatting again so no stray style id is left
# behind (matching rows 2-10, which carry no "s" attribute).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Rent",           "smws",  "2023-03-13", "1202.0", "Savings"),
    @("Rent",           "smw",   "2023-03-13", "1200.0", "Checkings"),
    @("Food",           "smwms", "2023-03-14", "150.0",  "Savings"),
    @("Transportation", "smwm",  "2023-03-14", "150.0",  "Checkings"),
    @("Clothing",       "swsd",  "2023-03-14", "450.0",  "Checkings")
)

$startRow = 11
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($col = 1; $col -le 5; $col++) {
        $text = [string]$values[$col - 1]
        $cell = $ws.Cells.Item($row, $col)

        # Force text so numeric/date-looking literals stay as strings.
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    }
}
